# Update the answer table: replace the division problems/answers in the
# five populated rows of the table with the new set of problems.
# The table always keeps 5 cells per row, so every change (including the
# ones that look like a cell insertion/removal in the raw OOXML diff)
# collapses to "this cell's text becomes that text" when walking the
# table row by row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rowUpdates = @{
    1  = @("59÷4=14, 3", "51÷5=10, 1", "38÷7=5, 3", "39÷8=4, 7", "84÷9=9, 3")
    5  = @("47÷9=5, 2", "55÷2=27, 1", "64÷5=12, 4", "21÷8=2, 5", "21÷9=2, 3")
    9  = @("38÷2=19, 0", "35÷9=3, 8", "85÷7=12, 1", "31÷9=3, 4", "56÷3=18, 2")
    13 = @("21÷2=10, 1", "83÷8=10, 3", "62÷6=10, 2", "77÷3=25, 2", "43÷3=14, 1")
    17 = @("77÷7=11, 0", "43÷9=4, 7", "65÷8=8, 1", "81÷8=10, 1", "92÷9=10, 2")
}

foreach ($rowIndex in $rowUpdates.Keys) {
    $values = $rowUpdates[$rowIndex]
    $row = $t.Rows.Item($rowIndex)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $cell = $row.Cells.Item($i + 1)
        $cell.Range.Text = $values[$i]
    }
}

Write-Host "table answers updated"
